$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: title & link update
$ws.Range("D16").Value = "[백준14916번, 그리디] 거스름돈 - Python"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/171"

# Row 46: title & link update
$ws.Range("D46").Value = "[CJ제일제당] 2021년 9월, 생물정보학(Bioinformatics 채용), BIO사업부문 신입R&D 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/415"

# Row 50: title & link update
$ws.Range("D50").Value = "Optuna"
$ws.Range("E50").Value = "http://incredible.egloos.com/7525213"

# Row 52: title update only
$ws.Range("D52").Value = "파이썬 관련 아르바이트 모집"
